$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.951.76"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").Value = "4.057.50"
$ws.Range("E3").Value = "  +3.30%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'523.26"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "'149.19"
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'0.739"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("D11").Value = "'0.0000341"
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("D12").Value = "'46.45"
$ws.Range("E12").Value = "  +9.67%  "
$ws.Range("E13").Value = "  +4.14%  "
$ws.Range("D14").Value = "4.704.65"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "4.058.59"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("D16").Value = "'21.52"
$ws.Range("E16").Value = "  +9.03%  "
$ws.Range("D17").Value = "'14.34"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").Value = "71.952.75"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").Value = "'442.85"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("D22").Value = "'3.54"
$ws.Range("E22").Value = "  +5.30%  "
$ws.Range("D23").Value = "'95.78"
$ws.Range("E23").Value = "  +9.00%  "
$ws.Range("D25").Value = "'14.41"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'4.10"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "'11.28"
$ws.Range("E27").Value = "  +5.18%  "
$ws.Range("D28").Value = "'37.35"
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'3.09"
$ws.Range("E30").Value = "  +8.78%  "
$ws.Range("D31").Value = "'702.73"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "'13.56"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.130"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "'6.94"
$ws.Range("E34").Value = "  +15.87%  "
$ws.Range("D35").Value = "'67.33"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D36").Value = "0.0₃0914"
$ws.Range("E36").Value = "  +8.62%  "
$ws.Range("D37").Value = "'0.448"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("D38").Value = "'41.09"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").Value = "'3.64"
$ws.Range("E39").Value = "  +22.77%  "
$ws.Range("D40").Value = "'0.155"
$ws.Range("E40").Value = "  +4.78%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("D44").Value = "'3.12"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").Value = "'3.53"
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("D48").Value = "'3.20"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "'0.000283"
$ws.Range("E49").Value = "  +24.91%  "
$ws.Range("D50").Value = "'9.17"
$ws.Range("E50").Value = "  +6.77%  "
$ws.Range("E51").Value = "  +1.69%  "
